$d = $word.ActiveDocument

# Insert a new paragraph before the current first paragraph ("Hello Git!")
$first = $d.Paragraphs(1).Range
$first.InsertParagraphBefore()

# The newly inserted (now empty) paragraph is the first paragraph; fill it in.
$greet = $d.Paragraphs(1).Range
$greet.Text = "Hallo Zusammen."
$greet.Font.Bold = $true
$greet.Font.Name = "Helvetica"
